$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 11 and row 12 for the columns that hold
# taxon-specific data (A, B, D, E, F, G, H, Q, R). The remaining columns
# (I, J, K, N, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY) are
# identical between the two rows, so no change is needed there.

$cols = @("A","B","D","E","F","G","H","Q","R")

foreach ($col in $cols) {
    $addr11 = "$col" + "11"
    $addr12 = "$col" + "12"

    $val11 = $ws.Range($addr11).Value2
    $val12 = $ws.Range($addr12).Value2

    $ws.Range($addr11).Value2 = $val12
    $ws.Range($addr12).Value2 = $val11
}
